$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.737.87"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "1.846.92"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "335.33"
$ws.Range("E5").Value = "  +0.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.007"
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4667"
$ws.Range("E7").Value = "  +0.81%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3848"
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.83"
$ws.Range("E9").Value = "  +1.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07903"
$ws.Range("E10").Value = "  +0.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9659"
$ws.Range("E11").Value = "  -3.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.26"
$ws.Range("E12").Value = "  -0.79%  "
$ws.Range("D13").Value = "1.849.71"
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.867"
$ws.Range("E14").Value = "  -1.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.120"
$ws.Range("E15").Value = "  -0.02%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.008"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "90.93"
$ws.Range("E17").Value = "  +2.97%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06614"
$ws.Range("E18").Value = "  -0.86%  "
$ws.Range("E19").Value = "  -0.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.20"
$ws.Range("E20").Value = "  +0.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.006"
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("D22").Value = "27.742.19"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.335"
$ws.Range("E23").Value = "  -0.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.80"
$ws.Range("E24").Value = "  -0.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.295"
$ws.Range("E25").Value = "  -0.70%  "
$ws.Range("D26").Value = "2.092.49"
$ws.Range("E26").Value = "  +1.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "158.87"
$ws.Range("E27").Value = "  +0.23%  "
$ws.Range("E28").Value = "  -0.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.061"
$ws.Range("E29").Value = "  -2.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.359"
$ws.Range("E30").Value = "  -0.76%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "118.45"
$ws.Range("E31").Value = "  -1.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09429"
$ws.Range("E32").Value = "  +0.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9369"
$ws.Range("E33").Value = "  -3.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.597"
$ws.Range("E34").Value = "  +0.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.247"
$ws.Range("E35").Value = "  -0.97%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.325"
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06020"
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02203"
$ws.Range("E38").Value = "  -0.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.202"
$ws.Range("E39").Value = "  -0.69%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.006"
$ws.Range("E40").Value = "  +0.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.154"
$ws.Range("E41").Value = "  -2.09%  "
$ws.Range("E42").Value = "  -1.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1843"
$ws.Range("E43").Value = "  -0.85%  "
$ws.Range("E44").Value = "  -2.71%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.300"
$ws.Range("E45").Value = "  +5.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.01"
$ws.Range("E46").Value = "  -1.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5436"
$ws.Range("E47").Value = "  -2.43%  "
$ws.Range("E48").Value = "  +1.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06826"
$ws.Range("E49").Value = "  +1.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "110.65"
$ws.Range("E50").Value = "  +0.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.007"
$ws.Range("E51").Value = "  -32.21%  "
